$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06416381063653058
$ws.Range("H2").Value = -0.9077691621289075
$ws.Range("I2").Value = 83.53132706545931
$ws.Range("G3").Value = 0.06417471588793254
$ws.Range("H3").Value = -6.18427187200777
$ws.Range("G4").Value = -0.04825625756548189
$ws.Range("H4").Value = -7.491629780685351
$ws.Range("G5").Value = -0.02343904928767914
$ws.Range("H5").Value = 13.4190838143734
$ws.Range("G6").Value = -0.09827839580595858
$ws.Range("H6").Value = 7.320553376493639
$ws.Range("G7").Value = -0.08679936669461225
$ws.Range("H7").Value = 5.003549060192296
$ws.Range("G8").Value = -0.3591808714796799
$ws.Range("H8").Value = 2.105086561854916
$ws.Range("G9").Value = -0.3691544364765446
$ws.Range("H9").Value = 5.367030519924915
$ws.Range("G10").Value = 0.03031715349703739
$ws.Range("H10").Value = 50.15183304104579
$ws.Range("G11").Value = 0.0407213349722756
$ws.Range("H11").Value = 79.43280510688707
$ws.Range("G12").Value = 0.2236250799009378
$ws.Range("H12").Value = 0.8466965848385432
$ws.Range("G13").Value = 0.2300624058376801
$ws.Range("H13").Value = 2.156064750365229
$ws.Range("G14").Value = -0.05454686181870474
$ws.Range("H14").Value = -29.54929048200705
$ws.Range("G15").Value = -0.05228807933285054
$ws.Range("H15").Value = -9.627126608247435
$ws.Range("G16").Value = 0.2141657570403719
$ws.Range("H16").Value = 0.749086328378764
$ws.Range("G17").Value = 0.2213484612587544
$ws.Range("H17").Value = 0.3629581708934411
$ws.Range("G18").Value = 0.07545180653851437
$ws.Range("H18").Value = 3.329988955619679
$ws.Range("G19").Value = 0.0863495800482556
$ws.Range("H19").Value = 14.61621677034472
$ws.Range("G20").Value = -0.0834552732223411
$ws.Range("H20").Value = -11.30044518973794
$ws.Range("G21").Value = -0.08606375840519892
$ws.Range("H21").Value = 0.5834351378092389
$ws.Range("G22").Value = 0.07415363741721973
$ws.Range("H22").Value = 0.8866989984645759
$ws.Range("G23").Value = 0.0697544579652344
$ws.Range("H23").Value = 2.082707735928478
$ws.Range("G24").Value = 0.05910846642883147
$ws.Range("H24").Value = -11.26471203893201
$ws.Range("G25").Value = 0.06329123637087247
$ws.Range("H25").Value = 15.53700993601712
$ws.Range("G26").Value = 0.119047007487078
$ws.Range("H26").Value = -0.2530617327840166
$ws.Range("G27").Value = 0.1247438616923543
$ws.Range("H27").Value = 9.55880983431245
$ws.Range("G28").Value = 0.1309681445486715
$ws.Range("H28").Value = 1.326669878909043
$ws.Range("G29").Value = 0.147674216408767
$ws.Range("H29").Value = -2.099037415724327
$ws.Range("G30").Value = 0.08856910946984774
$ws.Range("H30").Value = 5.055080872318008
$ws.Range("G31").Value = 0.09060846662767325
$ws.Range("H31").Value = 10.92010701105746
$ws.Range("G32").Value = 0.0554815329048789
$ws.Range("H32").Value = 3.976086773122307
$ws.Range("G33").Value = 0.05485697507196689
$ws.Range("H33").Value = -0.6986159075135691
$ws.Range("G34").Value = 0.01817063075344867
$ws.Range("H34").Value = 4.686150830192713
$ws.Range("G35").Value = 0.03180144817480884
$ws.Range("H35").Value = 88.16658776381431
$ws.Range("G36").Value = -0.02826839345418108
$ws.Range("H36").Value = 2.676151832599769
$ws.Range("G37").Value = -0.03502393934843685
$ws.Range("H37").Value = -5.291206992356573
$ws.Range("G38").Value = 0.07248532269285807
$ws.Range("H38").Value = -7.405721312489121
$ws.Range("G39").Value = 0.07725656677336626
$ws.Range("H39").Value = -0.6284176880287971
$ws.Range("G40").Value = 0.07121758393240943
$ws.Range("H40").Value = 7.564756924638334
$ws.Range("G41").Value = 0.06114691793410743
$ws.Range("H41").Value = -5.960861681135579
$ws.Range("G42").Value = 0.08612879324942375
$ws.Range("H42").Value = 10.71778059786821
$ws.Range("G43").Value = 0.08860783083255963
$ws.Range("H43").Value = 10.53251704166354
$ws.Range("G44").Value = 0.08415764427047931
$ws.Range("H44").Value = -4.634314846449596
$ws.Range("G45").Value = 0.09324547221413854
$ws.Range("H45").Value = 3.166470008684023
$ws.Range("G46").Value = 0.0004474921158187565
$ws.Range("H46").Value = 116.3543823860638
$ws.Range("G47").Value = 0.00200605051331723
$ws.Range("H47").Value = 4299.095458134407
$ws.Range("G48").Value = -0.09728128883374139
$ws.Range("H48").Value = -1.223810215725243
$ws.Range("G49").Value = -0.1055314843616371
$ws.Range("H49").Value = 3.691238266502343
$ws.Range("G50").Value = 0.1666912492339716
$ws.Range("H50").Value = -2.233464837360376
$ws.Range("G51").Value = 0.1750554836719478
$ws.Range("H51").Value = 3.073935930172054
$ws.Range("G52").Value = 0.07053348262663038
$ws.Range("H52").Value = -0.6037071411468897
$ws.Range("G53").Value = 0.06628213513103204
$ws.Range("H53").Value = 3.056819036836071
$ws.Range("G54").Value = -0.1369639202392226
$ws.Range("H54").Value = -7.168121981395307
$ws.Range("G55").Value = -0.1107729393551788
$ws.Range("H55").Value = 4.892393844214348
$ws.Range("G56").Value = 0.1858346636261513
$ws.Range("H56").Value = -2.203006578796351
$ws.Range("G57").Value = 0.2029276284079517
$ws.Range("H57").Value = 2.02070731356588
